# Insert two new data rows at row 292 (weekly update adds a new reporting
# date, pushing the rest of the "Femacal de La Calera - Piña" history down
# by two rows). Excel inherits formatting from the surrounding rows on
# insert, which already gives column D the date-number-format style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(292).Insert()
$ws.Rows.Item(292).Insert()

# New row 292: "Primera" quality entry for 2021-09-22 (serial 44461)
$ws.Cells.Item(292, 1).Value2 = 3
$ws.Cells.Item(292, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(292, 3).Value2 = "Coquimbo"
$ws.Cells.Item(292, 4).Value2 = 44461
$ws.Cells.Item(292, 5).Value2 = 5
$ws.Cells.Item(292, 6).Value2 = "Fruta"
$ws.Cells.Item(292, 7).Value2 = 100108
$ws.Cells.Item(292, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(292, 9).Value2 = 100108005
$ws.Cells.Item(292, 10).Value2 = "Piña"
$ws.Cells.Item(292, 11).Value2 = "Caramelo"
$ws.Cells.Item(292, 12).Value2 = "Primera"
$ws.Cells.Item(292, 13).Value2 = 108
$ws.Cells.Item(292, 14).Value2 = 20000
$ws.Cells.Item(292, 15).Value2 = 20000
$ws.Cells.Item(292, 16).Value2 = 20000
$ws.Cells.Item(292, 17).Value2 = "$/caja 12 unidades"
$ws.Cells.Item(292, 18).Value2 = "Ecuador"
$ws.Cells.Item(292, 19).Value2 = 1667
$ws.Cells.Item(292, 20).Value2 = 12

# New row 293: "Segunda" quality entry for 2021-09-22 (serial 44461)
$ws.Cells.Item(293, 1).Value2 = 3
$ws.Cells.Item(293, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(293, 3).Value2 = "Coquimbo"
$ws.Cells.Item(293, 4).Value2 = 44461
$ws.Cells.Item(293, 5).Value2 = 5
$ws.Cells.Item(293, 6).Value2 = "Fruta"
$ws.Cells.Item(293, 7).Value2 = 100108
$ws.Cells.Item(293, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(293, 9).Value2 = 100108005
$ws.Cells.Item(293, 10).Value2 = "Piña"
$ws.Cells.Item(293, 11).Value2 = "Caramelo"
$ws.Cells.Item(293, 12).Value2 = "Segunda"
$ws.Cells.Item(293, 13).Value2 = 162
$ws.Cells.Item(293, 14).Value2 = 20000
$ws.Cells.Item(293, 15).Value2 = 20000
$ws.Cells.Item(293, 16).Value2 = 20000
$ws.Cells.Item(293, 17).Value2 = "$/caja 14 unidades"
$ws.Cells.Item(293, 18).Value2 = "Ecuador"
$ws.Cells.Item(293, 19).Value2 = 1429
$ws.Cells.Item(293, 20).Value2 = 14
